# contactos_dante_propiedades.xlsx - "Add files via upload"
# New form responses were appended/edited upstream; this brings the
# "Contactos" sheet's logged rows (2-9) up to date and leaves the
# selection where the uploader last left it (C9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Fecha/Hora) and D (Teléfono) hold text that can look numeric
# (timestamps, phone numbers with leading zeros) -- force the Text number
# format first so values round-trip as strings instead of being coerced
# into numbers/dates.
$ws.Range("A2:A9").NumberFormat = "@"
$ws.Range("D2:D9").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "2025-10-07 20:56:03"
$ws.Range("B2").Value = "RAMONA"
$ws.Range("C2").Value = "DDDDDDD"
$ws.Range("D2").Value = "1151511579"

# Row 3
$ws.Range("A3").Value = "2025-10-07 21:03:33"
$ws.Range("B3").Value = "uf000"
$ws.Range("C3").Value = "GITHUB"
$ws.Range("D3").Value = "1136809319"

# Row 4
$ws.Range("A4").Value = "2025-10-07 21:03:57"
$ws.Range("B4").Value = "uf003"
$ws.Range("C4").Value = "propia"
$ws.Range("D4").Value = "01151511579"

# Row 5
$ws.Range("A5").Value = "2025-10-07 21:04:24"
$ws.Range("B5").Value = "uf004"
$ws.Range("C5").Value = "XXXXX"
$ws.Range("D5").Value = "0111551511579"

# Row 6
$ws.Range("A6").Value = "2025-10-07 21:04:49"
$ws.Range("B6").Value = "uf009"
$ws.Range("C6").Value = "GITHUB"
$ws.Range("D6").Value = "1136809319"

# Row 7
$ws.Range("A7").Value = "2025-10-07 21:05:09"
$ws.Range("B7").Value = "uf012"
$ws.Range("C7").Value = "GITHUB"
$ws.Range("D7").Value = "1151511579"

# Row 8
$ws.Range("A8").Value = "2025-10-08 10:56:18"
$ws.Range("B8").Value = "MARIA ROSA ARONA"
$ws.Range("C8").Value = "XXXXX"
$ws.Range("D8").Value = "1151511579"

# Row 9
$ws.Range("A9").Value = "2025-10-08 11:18:55"
$ws.Range("B9").Value = "ARTURO"
$ws.Range("C9").Value = "GITHUB"
$ws.Range("D9").Value = "0111551511579"

# Leave the selection on C9 (last edited cell), tab marked active, as
# recorded in the saved workbook view state.
[void]$ws.Range("C9").Select()
